$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text-safe cell updates (values Excel would not reinterpret as numbers)
$ws.Range('D2').Value = '70.712.62'
$ws.Range('E2').Value = '  +0.20%  '
$ws.Range('D3').Value = '3.533.92'
$ws.Range('E3').Value = '  -0.77%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('E5').Value = '  +3.87%  '
$ws.Range('E6').Value = '  +0.03%  '
$ws.Range('D7').Value = '3.530.71'
$ws.Range('E7').Value = '  -0.70%  '
$ws.Range('E8').Value = '  -0.77%  '
$ws.Range('E9').Value = '  -0.09%  '
$ws.Range('E10').Value = '  +1.34%  '
$ws.Range('E11').Value = '  -2.71%  '
$ws.Range('E12').Value = '  +0.03%  '
$ws.Range('E13').Value = '  -0.33%  '
$ws.Range('E14').Value = '  -0.53%  '
$ws.Range('D15').Value = '4.102.46'
$ws.Range('E15').Value = '  -0.81%  '
$ws.Range('E16').Value = '  +0.96%  '
$ws.Range('E17').Value = '  -0.73%  '
$ws.Range('D18').Value = '3.529.88'
$ws.Range('E18').Value = '  -0.73%  '
$ws.Range('D19').Value = '70.863.82'
$ws.Range('E19').Value = '  +0.34%  '
$ws.Range('E20').Value = '  +1.90%  '
$ws.Range('E21').Value = '  +1.99%  '
$ws.Range('E22').Value = '  +0.05%  '
$ws.Range('E23').Value = '  -1.60%  '
$ws.Range('E24').Value = '  -0.62%  '
$ws.Range('E25').Value = '  +0.88%  '
$ws.Range('E26').Value = '  -0.42%  '
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('E28').Value = '  -1.79%  '
$ws.Range('E29').Value = '  -0.53%  '
$ws.Range('E30').Value = '  +0.46%  '
$ws.Range('E31').Value = '  -1.85%  '
$ws.Range('E32').Value = '  -1.88%  '
$ws.Range('E33').Value = '  -0.33%  '
$ws.Range('E34').Value = '  -4.53%  '
$ws.Range('E35').Value = '  -7.66%  '
$ws.Range('E36').Value = '  +3.73%  '
$ws.Range('E37').Value = '  +0.71%  '
$ws.Range('E38').Value = '  -1.16%  '
$ws.Range('E41').Value = '  -6.18%  '
$ws.Range('E42').Value = '  +0.76%  '
$ws.Range('D43').Value = '3.349.66'
$ws.Range('E43').Value = '  -0.85%  '
$ws.Range('D44').Value = '0.0₃0728'
$ws.Range('E44').Value = '  +3.47%  '
$ws.Range('E45').Value = '  -2.65%  '
$ws.Range('E46').Value = '  -1.70%  '
$ws.Range('E47').Value = '  -2.56%  '
$ws.Range('E48').Value = '  -4.72%  '
$ws.Range('E49').Value = '  -0.18%  '
$ws.Range('E50').Value = '  +1.30%  '
$ws.Range('E51').Value = '  +7.62%  '

# Numeric-looking values: force text storage per cell (so the trailing zeros /
# exact digit strings survive), then restore the default "Normal" style so no
# cell formatting (s=...) is left behind, matching the source diff exactly.
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '621.35'
$ws.Range('D5').Style = "Normal"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '172.46'
$ws.Range('D6').Style = "Normal"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.608'
$ws.Range('D8').Style = "Normal"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '46.26'
$ws.Range('D13').Style = "Normal"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.0000276'
$ws.Range('D14').Style = "Normal"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '8.44'
$ws.Range('D16').Style = "Normal"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '607.28'
$ws.Range('D17').Style = "Normal"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '17.69'
$ws.Range('D21').Style = "Normal"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '15.60'
$ws.Range('D24').Style = "Normal"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '97.61'
$ws.Range('D25').Style = "Normal"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '33.67'
$ws.Range('D29').Style = "Normal"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '8.12'
$ws.Range('D31').Style = "Normal"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '6.80'
$ws.Range('D34').Style = "Normal"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '616.10'
$ws.Range('D35').Style = "Normal"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.0493'
$ws.Range('D36').Style = "Normal"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '10.84'
$ws.Range('D37').Style = "Normal"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '56.79'
$ws.Range('D39').Style = "Normal"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '3.40'
$ws.Range('D41').Style = "Normal"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.311'
$ws.Range('D45').Style = "Normal"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.89'
$ws.Range('D46').Style = "Normal"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '31.89'
$ws.Range('D47').Style = "Normal"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '134.07'
$ws.Range('D50').Style = "Normal"
